# The BOPIS_TestData sheet gains a new "Order#" value in H2 (Comments/Order#
# column for the BOPIS row), mirroring the order-number text values already
# present in the other sheets (e.g. "88012143", "88012146"). It must be
# written as text (not a number) so it matches the existing shared-string
# convention used for this column.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOPIS_TestData")

# Enter the value with a leading apostrophe so Excel stores it as text
# rather than re-interpreting the all-digit string as a number.
$ws.Range("H2").Value = "'88012153"

# Typing a quote-prefixed value mints a fresh (quotePrefix) cell style.
# Re-apply the formatting already used by the rest of row 2 (e.g. G2) so
# H2 keeps sharing that same style instead of the new one.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null
